# dataprovide for login failed due to null pointer exception
#
# The test data workbook had no data for the login page's DataProvider,
# which was causing a NullPointerException. This repurposes the blank
# "Sheet2" into a "LoginPage" sheet (keeping its sheetId), moves it in
# front of "contacts" so it becomes the active/first tab, and populates
# it with username/password test data plus a support email with a
# mailto hyperlink.

$wb = $excel.ActiveWorkbook

# Rename the existing blank "Sheet2" to "LoginPage" and move it to be
# the first (active) sheet, ahead of "contacts".
$loginSheet = $wb.Worksheets.Item("Sheet2")
$loginSheet.Name = "LoginPage"
$loginSheet.Move($wb.Worksheets.Item(1))

$ws = $wb.Worksheets.Item("LoginPage")

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Bad-credentials test case
$ws.Range("A2").Value = "&@#$#$@"
$ws.Range("B2").Value = "Auth@Bha!"

# Valid-looking account row (email reused as the bad-password test value)
$ws.Range("A3").Value = "pussharma@gmail.com"
$ws.Range("B3").Value = "&@#$#$@"

# Excel auto-hyperlinks an email address typed into a cell
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:pussharma@gmail.com")

# Widen column A to fit the longest value ("pussharma@gmail.com")
$ws.Range("A1:A3").EntireColumn.AutoFit()

# Restore the selection that was active on this sheet
$ws.Range("B6").Select()
